$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.769.40'
$ws.Range("E2").Value = '  -0.41%  '

$ws.Range("E3").Value = '  -1.38%  '

$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  -0.33%  '

$ws.Range("D5").Value = '314.45'
$ws.Range("E5").Value = '  -1.07%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("E7").Value = '  -2.65%  '

$ws.Range("D8").Value = '0.3960'
$ws.Range("E8").Value = '  -2.85%  '

$ws.Range("D9").Value = '52.07'
$ws.Range("E9").Value = '  -2.77%  '

$ws.Range("D10").Value = '1.000'
$ws.Range("E10").Value = '  -0.28%  '

$ws.Range("D11").Value = '1.404'
$ws.Range("E11").Value = '  -5.29%  '

$ws.Range("D12").Value = '0.08661'
$ws.Range("E12").Value = '  -1.96%  '

$ws.Range("D13").Value = '25.34'
$ws.Range("E13").Value = '  -4.21%  '

$ws.Range("D14").Value = '7.333'
$ws.Range("E14").Value = '  -2.13%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '7.783'
$ws.Range("E15").Value = '  -4.53%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.00001320'
$ws.Range("E16").Value = '  -2.95%  '

$ws.Range("D17").Value = '1.684.62'
$ws.Range("E17").Value = '  -1.39%  '

$ws.Range("D18").Value = '93.70'
$ws.Range("E18").Value = '  -3.50%  '

$ws.Range("D19").Value = '0.07070'
$ws.Range("E19").Value = '  -1.47%  '

$ws.Range("D20").Value = '20.34'
$ws.Range("E20").Value = '  -4.26%  '

$ws.Range("D21").Value = '7.082'
$ws.Range("E21").Value = '  -2.87%  '

$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").Value = '13.95'
$ws.Range("E23").Value = '  -3.07%  '

$ws.Range("D24").Value = '24.748.27'
$ws.Range("E24").Value = '  -0.48%  '

$ws.Range("D25").Value = '2.348'
$ws.Range("E25").Value = '  +0.92%  '

$ws.Range("D26").Value = '2.802'
$ws.Range("E26").Value = '  -4.28%  '

$ws.Range("D27").Value = '23.47'

$ws.Range("D28").Value = '162.57'
$ws.Range("E28").Value = '  -2.62%  '

$ws.Range("D29").Value = '5.843'
$ws.Range("E29").Value = '  -8.06%  '

$ws.Range("D30").Value = '147.58'
$ws.Range("E30").Value = '  +1.11%  '

$ws.Range("D31").Value = '7.904'
$ws.Range("E31").Value = '  -6.24%  '

$ws.Range("D32").Value = '2.408'
$ws.Range("E32").Value = '  +7.84%  '

$ws.Range("D33").Value = '1.943.64'
$ws.Range("E33").Value = '  +2.64%  '

$ws.Range("D34").Value = '0.08433'
$ws.Range("E34").Value = '  -4.65%  '

$ws.Range("D35").Value = '0.03062'
$ws.Range("E35").Value = '  -4.29%  '

$ws.Range("D36").Value = '6.977'
$ws.Range("E36").Value = '  -3.80%  '

$ws.Range("D37").Value = '1.003'
$ws.Range("E37").Value = '  -3.00%  '

$ws.Range("D38").Value = '0.2820'
$ws.Range("E38").Value = '  -2.16%  '

$ws.Range("D39").Value = '0.09483'
$ws.Range("E39").Value = '  +1.31%  '

$ws.Range("D40").Value = '10.59'
$ws.Range("E40").Value = '  -2.83%  '

$ws.Range("D41").Value = '1.510'
$ws.Range("E41").Value = '  +2.60%  '

$ws.Range("D42").Value = '0.7962'
$ws.Range("E42").Value = '  -5.81%  '

$ws.Range("D43").Value = '13.60'
$ws.Range("E43").Value = '  -4.11%  '

$ws.Range("E44").Value = '  -5.26%  '

$ws.Range("D45").Value = '0.7169'
$ws.Range("E45").Value = '  -3.64%  '

$ws.Range("D46").Value = '2.573'
$ws.Range("E46").Value = '  -5.14%  '

$ws.Range("D47").Value = '4.180'
$ws.Range("E47").Value = '  -1.47%  '

$ws.Range("D48").Value = '0.08712'
$ws.Range("E48").Value = '  +4.22%  '

$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D49").Value = '1.000'
$ws.Range("E49").Value = '  +0.22%  '

$ws.Range("B50").Value = 'Flow'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D50").Value = '1.352'
$ws.Range("E50").Value = '  -3.63%  '

$ws.Range("D51").Value = '138.25'
$ws.Range("E51").Value = '  -2.58%  '
